$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'white athletic leggings mens'
    'supreme basketball tights'
    'snowmobiling knee pads'
    'photographer knee pads'
    'raymens leggings'
    'training basketball youth'
    'baleaf men compression pants'
    'winter leggings men'
    'kids pants with knee pads'
    'eclipse knee pads'
    'pilates knee pads'
    'fitted mens tights'
    'mens knee pads bmx'
    'canoe knee pad'
    'graduated compression pants'
    'graduated compression leggings'
    'soft knee pads'
    'knee pads nike'
    'knee pad toddler'
    'knee pad yoga'
    'knee pads dance'
    'knee pads mizuno'
    'elbow knee pads'
    'nee pads basketball'
    'compression pants set men'
    'mens basketball snap pants'
    'nike compression pants for youth'
    'mens compression tights cold weather'
    'mens under armour basketball tights'
    'mens tights pockets'
    'rollerblades knee pads'
    'knee pads 3xl'
    'knee pads 8'
    'knee pad hard'
    'men leggings fleece'
    'teflex knee pads'
    'man leggings thermal'
    'sailing knee pad'
    'knee pads downhill'
    'knee pads airsoft'
    'knee pads army'
    'knee pads enduro'
    'knee pads bike'
    'knee pads caving'
    'knee pads for women'
    'knee pads gloves'
    'knee pads green'
    'knee pads kali'
    'knee pads kuangmi'
    'knee pads longboard'
    'knee pads neoprene'
    'knee pads ocp'
    'knee pads orange'
    'knee pads over pants'
    'knee pads plastic'
    'knee pads purple'
    'knee pads razor'
    'knee pads red'
    'knee pads rollerblading'
    'knee pads sailing'
    'knee pads scooter'
    'knee pads set'
    'knee pads shooting'
    'knee pads swat'
    'knee pads teen'
    'knee pads tsg'
    'knee pads usmc'
    'knee pads viper'
    'knee pads white'
    'knee pads yellow'
    'mens leggings xs'
    'xtextile compression pants men'
    'crx men''s tights'
    'yoga capri pants'
    'dodoing kneepads'
    'cavaliers basketball leggings'
    'basketball knee pads kids'
    'knee pads for toddlers'
    'knee pads skating'
    'knee pads skateboarding'
    'basketball knee pads kids boys'
    'knee pads for dance'
    'knee pads rollerblade'
    'knee pads tan'
    'knee pad dancer'
    'knee pads adidas'
    'knee pads basketball mcdavid'
    'knee pads dancing'
    'knee pads dodgeball'
    'knee pads pair'
    'knee pads longboarding'
    'knee pads nba'
    'knee pads pole'
    'knee pad and elbow pads'
    'knee pad adidas'
    'knee pad asics'
    'knee pad for kids'
    'knee pad military'
    'knee pad mma'
    'knee pad mizuno'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
